$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)

# --- Sheet1 ("Image CLEF2013 natural leaf"): zoom 182 -> 120 ---
$sheet1.Activate()
$excel.ActiveWindow.Zoom = 120

# --- Sheet2 ("CLEF2013 uniform leaf"): the bulk of the edits ---
$sheet2.Activate()

# Existing L2 just needs the wrap-text style applied (value/string unchanged)
$sheet2.Range("L2").WrapText = $true

# Row 10: update existing test-accuracy value, add a new training-accuracy (C10)
# value, and swap the comment string + wrap it.
$sheet2.Range("B10").Value = 0.85599999999999998
$sheet2.Range("C10").Value = 0.40949999999999998
$sheet2.Range("L10").Value = "Set threshold to be >= 20 samples per species; K=1000. No resize before SIFT."
$sheet2.Range("L10").WrapText = $true

# Row 11 (new): SIFT_BOF / SVM linear run with K=3000
$sheet2.Range("A11").Value = "SIFT_BOF"
$sheet2.Range("B11").Value = 0.86629999999999996
$sheet2.Range("C11").Value = 0.45900000000000002
$sheet2.Range("D11").Value = 109
$sheet2.Range("E11").Value = 9607
$sheet2.Range("F11").Value = 66
$sheet2.Range("G11").Value = 1194
$sheet2.Range("H11").Value = "N/A"
$sheet2.Range("I11").Value = "no"
$sheet2.Range("J11").Value = "N/A"
$sheet2.Range("K11").Value = "SVM linear"
$sheet2.Range("L11").Value = "Set threshold to be >= 20 samples per species; K=3000. No resize before SIFT."
$sheet2.Range("L11").WrapText = $true

# Row 4 (new): TL / ResNet50 run with histogram normalization + K-means background removal
$sheet2.Range("A4").Value = "TL"
$sheet2.Range("B4").Value = 0.99860000000000004
$sheet2.Range("C4").Value = 5678
$sheet2.Range("D4").Value = 109
$sheet2.Range("E4").Value = 28821
$sheet2.Range("F4").Value = 66
$sheet2.Range("G4").Value = 1194
$sheet2.Range("H4").Value = "no"
$sheet2.Range("I4").Value = "factor=3.`nrotation=10.`nshift=0.1"
$sheet2.Range("I4").WrapText = $true
$sheet2.Range("J4").Value = "ResNet50"
$sheet2.Range("K4").Value = "SVM linear"
$sheet2.Range("L4").Value = "Set threshold to be >= 20 samples per species. Apply histogram normalization and then use K-means to remove background color "
$sheet2.Range("L4").WrapText = $true
$sheet2.Rows.Item(4).RowHeight = 48

# Column widths: new column I, and give column L (already wrap-texted) its
# dedicated style column.
$sheet2.Columns.Item(9).ColumnWidth = 17.83
$sheet2.Columns.Item(12).ColumnWidth = 67.83

# View state: scroll so column E is left-most, and land the selection on M15.
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$sheet2.Range("M15").Select()
